# Update row 9 (Ano 2025) figures in the faturamento anual sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = 3912398.62
$ws.Range("C9").Value = 619711.24
$ws.Range("D9").Value = 4532109.86
$ws.Range("E9").Value = 13.67379121740884
$ws.Range("F9").Value = 86.32620878259117
$ws.Range("G9").Value = -40.10799867669782
$ws.Range("H9").Value = -29.3474887951604
$ws.Range("I9").Value = 39761
$ws.Range("J9").Value = 1703
$ws.Range("K9").Value = 41464
$ws.Range("L9").Value = 28699
$ws.Range("M9").Value = 157.9187379351197
$ws.Range("N9").Value = 7.814338834429235
